$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$cell = $tbl.Cell(2,1)
$tf = $cell.Shape.TextFrame
$tr = $tf.TextRange
# Step 1: blank out run3 "(Hank)" (positions 4-9)
$r3 = $tr.Characters(4,6)
$r3.Text = ""
Write-Host "After step1:" $tr.Text
# Step2: blank out whats-now-run2 (translated "翰", at position...?)
$full = $tr.Text
Write-Host "Full after step1 (repeated):" $full
